$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.907.75"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.503.48"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.17"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.38"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "3.502.44"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.192"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.28"
$ws.Range("E11").Value = "  +7.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.582"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.12"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000274"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "4.064.78"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.28"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "611.20"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "69.939.78"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "3.495.93"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.873"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.09"
$ws.Range("E23").Value = "  -18.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.55"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.92"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.71"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.56"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.24"
$ws.Range("E29").Value = "  +4.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.10"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.96"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "652.60"
$ws.Range("E33").Value = "  +14.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.93"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.58"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0999"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.70"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").Value = "  +8.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.36"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").Value = "3.322.47"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.92"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.30"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "0.0₃0689"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.55"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.31"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("E51").Value = "  -0.04%  "
